$wb = $excel.ActiveWorkbook

# --- Sheet "constants": remove the duplicated Bulgaria-specific "improve_dst"
# parameter block (rows 36-40), which duplicated the generic "improve_dst"
# parameters already present in rows 31-35. Deleting the rows shifts
# everything below up by 5 rows (old 41-55 -> new 36-50).
$ws1 = $wb.Worksheets.Item("constants")
$ws1.Range("A36:G40").EntireRow.Delete()

# Update the view/selection on the constants sheet.
$ws1.Application.ActiveWindow.ScrollRow = 31
$ws1.Range("A39").Select()

# --- Sheet "time_variants": the row that used to read
# "int_perc_bulgaria_improve_dst" is renamed to the generic
# "int_perc_improve_dst" (now that the Bulgaria-specific duplicate set of
# parameters has been removed from "constants").
$ws2 = $wb.Worksheets.Item("time_variants")
$ws2.Range("A6").Value2 = "int_perc_improve_dst"

# A handful of scenario percentages on this sheet were bumped from 99 to 100.
$ws2.Range("T6").Value2 = 100
$ws2.Range("S7").Value2 = 100
$ws2.Range("M8").Value2 = 100
$ws2.Range("R8").Value2 = 100
$ws2.Range("S8").Value2 = 100

# V10 had a stray 0 value cleared out (style/formatting left untouched).
$ws2.Range("V10").ClearContents()

# Update the view/selection on the time_variants sheet.
$ws2.Select()
$ws2.Application.ActiveWindow.Zoom = 115
$ws2.Application.ActiveWindow.FreezePanes = $false
$ws2.Range("J2").Select()
$ws2.Application.ActiveWindow.FreezePanes = $true
$ws2.Range("A9").Select()
